# Commit: "add some new equips. remove some useless module"
#
# This updates the People sheet (人物):
#  - Column O (BattleMap) for rows 14-42 changes from "oneline" to "doublec"
#    (dropping the extra "new equip" / module variants down to the
#    "doublec" battle map setting for those rows)
#  - The People sheet becomes the active/selected sheet, with cell O16
#    selected as the active cell in the bottom-right frozen pane
#  - Whatever sheet was previously active (怪物) loses its "active" flag
#    automatically once a different sheet is activated

$wb = $excel.ActiveWorkbook

$wsPeople = $wb.Worksheets.Item("人物")

# Update the BattleMap column (O) for rows 14 through 42
$wsPeople.Range("O14:O42").Value = "doublec"

# Make the People sheet the active sheet, and select O16 as the active cell
$wsPeople.Activate() | Out-Null
$wsPeople.Range("O16").Select() | Out-Null
